$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) for rows 2-9 with the newly regenerated values
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 11
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 10
$ws.Range("G6").Value = 10
$ws.Range("G7").Value = 7
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 4
